$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Distance1"
$ws.Range("B1").Value = "Sensor 1"
$ws.Range("C1").Value = "Distance3"
$ws.Range("D1").Value = "Sensor 3"

$ws.Range("F3").Select() | Out-Null
